$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: write column B (text) in the exact order new distinct strings were
# first introduced in the source workbook, so the shared-string table the
# engine (re)builds on save lines up with the target uniqueCount ordering.
$bTexts = @(
    @(93, 'Bordo(A-B) yurtlarda son iki haftadır su kıtlığı çekmeye başladık. Özellikle geceleri asla içme suyu olmuyor. Lütfen geldiği zaman suyu odalarınızda stoklamayın.'),
    @(94, 'Erasmu stajı ile ilgili bir sorum olacaktı. İngilizce okuyorum ama hazırlığı direkt giriş yeterlilik sınavında atlamıştım ve exit sınavına hiç girmedim. Başvuru yaparken sistem yabancı dil belgesini zorunlu olarak istiyor. Sınav ne zaman oluyor nasıl başvuruluyor bilgisi olan var mı?'),
    @(95, 'Yurt oda değişikliği oluyor mu hiçbir fikrim yok.'),
    @(96, 'Merhaba developer student diye bir klüp var değil mi? Gördüm diye hatırlıyorum ama bulamıyorum.'),
    @(97, 'Ders programında "yes" yazanlar online mı?'),
    @(98, 'Merhaba, yurt iptali için son gün ne zaman bilen var mı?'),
    @(99, 'Merhaba yurt ödemeleri ile alakalı işlem yapabilen var mı? Veya ödeme için sistemin açılacağı bir tarih varsa yazabilir misiniz?'),
    @(100, 'Yurt iptal işlemi için en son gün ne zaman?'),
    @(101, 'Dönem dondurma dilekçemi aktardım fakat FYK kurulu toplanacak karar bildirilecek diyor yarın mali yükümlülükler için son gün problem olur mu bilgisi olan var mı?'),
    @(102, 'Yurt için gerekli belgeler neler ve ne zaman teslim ediyoruz?'),
    @(103, 'AKTS ücreti ödenirken nasıl bir yol izleniyor? Ödeme yaparken ne yapmak lazım?'),
    @(104, 'Yurt iptali Pazartesi mesai bitimine kadarmış, hangi derslerin online olduğu da Pazartesi açıklanacak. Pazartesi tüm derslerim online olursa iptal etmeyi düşünüyorum ancak dersler akşam mı açıklanır sabah mı bilgisi olan var mı?'),
    @(105, 'Yurt iptali veya birine devretmek için ne yapmam gerekiyor?'),
    @(106, 'Shuttle saatleri ne zaman belli olacak?'),
    @(107, 'İngilizce inşaat mühendisliği koşullu derslerin koşuluna falan nasıl ulaşabilirim sitesinde bulamadım da.'),
    @(108, 'Okul ücretini AKTS üzerinden yatırırken nasıl bir yol izleniyor? Parayı yatırırken açıklama kısmına mı yazmak gerekiyor?'),
    @(109, 'Online post sistemi ile bu senenin ödemesini yaptım ama e-campus''de hala Holds%Warnings bölümünde ödeme uyarısı var aynı problemi yaşayan başkası var mı?'),
    @(110, 'Bursumuz keslince 27 kredi alabiliyoruz diye duydum doğru mu acaba?'),
    @(111, 'Yurt fiyatları neden bu kadar yüksek öğrenci konseyi bununla alakalı bir çalışma yapıyor mu?'),
    @(112, 'Çalışma burslu öğrenciler için yurt ayarlanıyor mu?'),
    @(113, 'Yurt başvurumu nasıl iptal edebilirim? Mail attım geri dönüş olmadı.'),
    @(114, 'Bu dönem devamsızlık olacak mı?'),
    @(115, 'Yurt tercihleri neden pandemi kararları alınmadan önce yapıldı? Belki de bu kararlardan ötürü yurt tutmayacaklar olacaktı. Yurtlar müdürlüğü bu kararların ardından iyi niyet gösterip iade isteyenlere tam iade yapacak mı? Yoksa bu zamanlamayı öğrenciler ticari bi strateji olarak mı değerlendirmeli?'),
    @(116, 'Hibrit modelinde şehir dışında yaşayan öğrenciler için çevrimiçi ders önceliği verilir mi acaba?'),
    @(117, 'Erasmus hibe konusunda bilgisi olan var mı? Kabul mektubu teslim süresinden 10 gün geçti ve hala okul bir şey açıklamadı.'),
    @(118, 'Yurt sonuçları ne zaman açıklanacak belli mi?'),
    @(119, '2017''de kayıt oldum, tek dersim kaldı onun için akts olarak mı yoksa kredi olarak mı ödemesini yapacağım?'),
    @(120, 'Okulda son dönemim ve dönem içinde staj yapacağım bunun için stajı ders olarak eklemeli miyim?'),
    @(121, 'Staj yapacağımız şirkette ilgili mühendislikten en az 2 çalışan olması lazım gibi bir kural var mı?'),
    @(122, 'Ben okuldan mezun oldum ama bana okul kep göndermiş. Bunu iade edecek miyim? Ben 2 sene önce mezun oldum niye kep gönderirler ki?'),
    @(123, 'Geri ödemeli bursun başvuru süreci hakkında bilgi sahibi olan var mı? Yönergede nasıl başvurulacağına dair detayları bulamadım.'),
    @(124, 'Tek ders sınavları için okula dilekçe yazdım fakat hala geri dönüş olmadı. Benime aynı durumda olan var mı?'),
    @(127, 'Arkadaşlar zorunlu staj formlarının okul tarafından onaylanması tahminen kaç gün sürer?'),
    @(128, 'Erasmus''a öğrenci seçilirken aynı okula tercih yapmış okuldaki bütün öğrenciler arasından mı seçiliyor yoksa aynı bölüm içindeki öğrencilerden mi seçiliyor? Kısacası, farklı bölümden öğrenciler birbirine rakip olabiliyor mu?'),
    @(129, 'Maslak kampüsten gelen öğrencilere indirim yapın onlar yurtları doldursun. Yurtlarda yer kalmayacak evlere akın edeceğiz ev kiraları 2 katına çıkar artık.'),
    @(130, 'Maslak''ta okuyan öğrencilerin bizden ne farkları var? Ya da şöyle sorayım bizim onlardan ne eksiğimiz var? Onlar indirimden yararlanırken biz neden yararlanamıyoruz?'),
    @(131, 'Başka okuldan ders alacam fakat bizim okulla uyumlu olan hangi üniversite var bilgisi olan var mı?'),
    @(132, 'Staj yaparken yaz okulu yapabilir miyiz?'),
    @(133, 'Bizim okuldan aldığım ders ile başka bir okuldan aldığım dersin 3 saati çakışıyor. Bu konu hakkında ne yapabilirim?'),
    @(134, 'Son dersimi başka bir okuldan almak için okula dilekçe verdim, kabul ettiler. Ama ders aldığım okulda kaydın bitmesine 1 gün kala kota dolmayacağı görünüyor. Başka bir okuldan alabilir miyim?'),
    @(135, 'Erasmus başvurusunda yazılı dil puanı ve sözlü mülakattan barajı geçince başka bir eleme tarzı bir olay var mı? Yoksa barajı geçen herkes yerleşebiliyor mu?'),
    @(136, 'Havale/Eft ile okula ödeme yapabiliyor muyuz?'),
    @(137, 'Yaz okulu için akademik takvimde son ödeme günü 9 Temmuz yazıyor, doğru mu?'),
    @(138, 'Bir staj yeri bulduktan sonra bunu okula nasıl bildiriyoruz? Nereye nasıl yazmak gerekir? '),
    @(139, 'Bölümüm İngilizce başka bir üniversitenin yaz okulundan istediğim dersin içeriği aynı ise Türkçe''sini alabilir miyim?'),
    @(140, '2015 girişliyim hazırlığı atladım bölüme başladım. Bu sene 6. senem bitti bölümde. 7 sene sonunda okuldan atılıyor muyum? Bursu kesilen ve son sınıf öğrencisi olan önkoşul olayına takılıyor mu?'),
    @(141, 'Yüz yüze eğitim için geleceğimiz haftayı değiştirebiliyor muyuz?'),
    @(142, '1. aşıyı oldum ama 2. aşının zamanı gelmedi henüz, pazartesi okula gelirken ocr testi vermeme gerek var mı?'),
    @(143, 'Bu dönem bölümü dondurdum. Kyk''dan kredi alıyordum. Kredi kesilir mi?'),
    @(144, 'Ders programım yeni onaylandı Blackboard''da ne zaman işlenir bilgisi olan var mı?'),
    @(145, 'Derslerde 2 çakışmanın kabul edilmesi için hiç mi çare yok bilgisi olan var mı?'),
    @(146, 'Mezun olmak için 4 dersim kaldı. DD derslerimi almak zorunda mıyım?'),
    @(147, 'Ben DGS sınavı ile geldim. Staj olaylarını bütün dersleri verdikten sonra yapmak gibi bir durum söz konusu oluyor mu?'),
    @(148, 'Ödeme için 15''i son gün ve dönemi dondurmak istiyorum. Yine de dilekçeyi vermeden önce ödeme yapmam gerekiyor mu?'),
    @(149, 'Bilgisayar mühendisliğinde 3. sınıf sayılabilmek için toplamda kaç krediyi vermiş olmamız gerekiyor?'),
    @(150, 'Güz dönemi bi dersten çekilmiştim. Şimdi bahar dönemi tekrar açılıyor. O dersi almak zorunda mıyım, sistem açıldığı ilk dönem alma koşulu koyuyor mu?'),
    @(151, 'Çap yapmak için dönem ortalamamız 2.5 ve F olmamalı değil mi?'),
    @(152, 'Tek ders sınavı tarihleri vs. hakkında bilgisi olan var mı?'),
    @(153, 'Ben final sınavlarının çoğuna yüz yüze girecek kişilerden biriyim ve vakaların böyle arttığı bir dönemde gelmek istemiyorum. Uzunca bir süredir korkumdan çok gerekmedikçe evden çıkmıyorum. Okula gelirken 3 vasıta değiştiriyorum ve hepsi de kalabalık oluyor. Bu konu hakkında bir şey yapılabilir mi?'),
    @(154, 'Okul haftaya eğitimin yüz yüze devam edip etmeyeceğiyle ilgili net bir açıklama yaptı mı? Şehir dışından gelen öğrenciler için bunun önceden bildirilmiş olması önemli. Yurda gelip de 2 gün sonra geri dönmek istemiyorum.'),
    @(155, 'Metrogarden''dan kalkan servisler AVM''nin otoparkında mı oluyor?'),
    @(156, 'Bir bölümün ingilizcesinin türkçesine nasıl geçiliyor?'),
    @(157, 'Hazırlıkta kaç saat devamsızlık hakkımız var?'),
    @(158, 'Mezun olmama 4 ders kaldı fakat 3 tanesi bu dönem açıldı. Okul benden dönem parası istiyor, mali işler 3 ders parası verebilirsin diyor anlamadım gitti.'),
    @(159, 'Add-drop haftasından sonra kayıt yaptıran var mı? Son tarihten sonra kayıt yaptırsam ne olur?'),
    @(160, 'Bu dönem ücretimi yatırdım ama şuan dondurmak istiyorum. Dondurmak için ücret veriliyor mu? Verilmiyorsa ödediğim ücreti geri alabiliyor muyum?'),
    @(161, 'Staj sürecinin işleyişi hakkında bilgi verebilecek olan var mı?'),
    @(162, 'Ders programından ders dün ders çıkarıldı ve programlarımız onaylanmışken oldu bu. Yeniden sistem açtırıp ders seçmek zorunda kaldık. Bu durumda birçok derste kota sorunu ve çakışma ile karşılaştık.'),
    @(163, 'Benim son senem nu sene mezun olmam gerekiyor fakat pandemiden ötürü staj yapamadım. Yaz tatilinde yapsam ne kadar ücret ödemem gerekiyor?'),
    @(164, 'Ders seçimi gece mi açılıyor sabah 9''da mı?'),
    @(165, 'ÇAP başvuruları ne zaman oluyor?'),
    @(166, 'ÇAP için minimum not ortalamamızın kaç olması gerekiyor?'),
    @(167, 'ÇAP yapınca YKS bursu kesiliyor mu?'),
    @(168, 'Döneme irregular olarak başlayabiliyor muyum?'),
    @(169, 'Disiplin soruşturması açılırsa burs kesilmesi oluyor mu?'),
    @(125, 'Bir derste kotaya ihtiyacım var. E-campus üzerinden istek yolluyorum. Danışmanım kabul ediyor fakat kota yöneticisi hala görmedi. Kime mail atmam gerekir, nereye ulaşmam lazım?'),
    @(126, 'Şu an ders ekleme bırakma haftasındayız fakat E-campus''te withdraw kapalı yazıyor. Bilgisi olan var mı?')
)
foreach ($pair in $bTexts) {
    $ws.Cells.Item($pair[0], 2).Value = $pair[1]
}

# Step 2: write columns A (id) and C (subject) for each new row (order
# doesn't matter here -- ids are numeric, subjects already exist in the
# shared-string table).
$rowsData = @(
    @(93, 91, 'Yurtlar'),
    @(94, 92, 'Öğrenci İşleri'),
    @(95, 93, 'Yurtlar'),
    @(96, 94, 'Sağlık, Kültür ve Spor'),
    @(97, 95, 'Öğrenci İşleri'),
    @(98, 96, 'Yurtlar'),
    @(99, 97, 'Mali İşler'),
    @(100, 98, 'Yurtlar'),
    @(101, 99, 'Öğrenci İşleri'),
    @(102, 100, 'Yurtlar'),
    @(103, 101, 'Mali İşler'),
    @(104, 102, 'Öğrenci İşleri'),
    @(105, 103, 'Yurtlar'),
    @(106, 104, 'Servis'),
    @(107, 105, 'Öğrenci İşleri'),
    @(108, 106, 'Mali İşler'),
    @(109, 107, 'Mali İşler'),
    @(110, 108, 'Öğrenci İşleri'),
    @(111, 109, 'Yurtlar'),
    @(112, 110, 'Öğrenci İşleri'),
    @(113, 111, 'Yurtlar'),
    @(114, 112, 'Öğrenci İşleri'),
    @(115, 113, 'Yurtlar'),
    @(116, 114, 'Öğrenci İşleri'),
    @(117, 115, 'Öğrenci İşleri'),
    @(118, 116, 'Yurtlar'),
    @(119, 117, 'Öğrenci İşleri'),
    @(120, 118, 'Öğrenci İşleri'),
    @(121, 119, 'Öğrenci İşleri'),
    @(122, 120, 'Öğrenci İşleri'),
    @(123, 121, 'Mali İşler'),
    @(124, 122, 'Öğrenci İşleri'),
    @(125, 123, 'Öğrenci İşleri'),
    @(126, 124, 'Öğrenci İşleri'),
    @(127, 125, 'Öğrenci İşleri'),
    @(128, 126, 'Öğrenci İşleri'),
    @(129, 127, 'Dekanlık'),
    @(130, 128, 'Dekanlık'),
    @(131, 129, 'Öğrenci İşleri'),
    @(132, 130, 'Öğrenci İşleri'),
    @(133, 131, 'Öğrenci İşleri'),
    @(134, 132, 'Öğrenci İşleri'),
    @(135, 133, 'Öğrenci İşleri'),
    @(136, 134, 'Mali İşler'),
    @(137, 135, 'Dekanlık'),
    @(138, 136, 'Öğrenci İşleri'),
    @(139, 137, 'Öğrenci İşleri'),
    @(140, 138, 'Öğrenci İşleri'),
    @(141, 139, 'Öğrenci İşleri'),
    @(142, 140, 'Sağlık, Kültür ve Spor'),
    @(143, 141, 'Öğrenci İşleri'),
    @(144, 142, 'Uzem'),
    @(145, 143, 'Öğrenci İşleri'),
    @(146, 144, 'Öğrenci İşleri'),
    @(147, 145, 'Öğrenci İşleri'),
    @(148, 146, 'Mali İşler'),
    @(149, 147, 'Öğrenci İşleri'),
    @(150, 148, 'Öğrenci İşleri'),
    @(151, 149, 'Öğrenci İşleri'),
    @(152, 150, 'Öğrenci İşleri'),
    @(153, 151, 'Dekanlık'),
    @(154, 152, 'Dekanlık'),
    @(155, 153, 'Servis'),
    @(156, 154, 'Öğrenci İşleri'),
    @(157, 155, 'Öğrenci İşleri'),
    @(158, 156, 'Mali İşler'),
    @(159, 157, 'Öğrenci İşleri'),
    @(160, 158, 'Mali İşler'),
    @(161, 159, 'Öğrenci İşleri'),
    @(162, 160, 'Öğrenci İşleri'),
    @(163, 161, 'Öğrenci İşleri'),
    @(164, 162, 'Öğrenci İşleri'),
    @(165, 163, 'Öğrenci İşleri'),
    @(166, 164, 'Öğrenci İşleri'),
    @(167, 165, 'Öğrenci İşleri'),
    @(168, 166, 'Öğrenci İşleri'),
    @(169, 167, 'Dekanlık')
)
foreach ($item in $rowsData) {
    $ws.Cells.Item($item[0], 1).Value = $item[1]
    $ws.Cells.Item($item[0], 3).Value = $item[2]
}

$ws.Range("F149").Select() | Out-Null